$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had years 2015..2020 in columns B..G (row 1), with
# summary stats (count/mean/std/min/25%/50%/75%/max) below.
# The new version prepends years 2010..2014 (5 more columns) before the
# existing data, shifting the old B..G columns to G..L, and re-computes a
# few statistics whose underlying values changed slightly.

# 1) Insert 5 new blank columns before column B. This pushes the existing
#    B:G data to G:L.
$ws.Range("B1:F1").EntireColumn.Insert()

# 2) The new header cells (row 1) need the same bold/bordered style as the
#    other year header cells. Copy formatting (only) from an existing
#    header cell onto the new ones, then fill in the year values.
$ws.Range("H1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B1").Value = 2010
$ws.Range("C1").Value = 2011
$ws.Range("D1").Value = 2012
$ws.Range("E1").Value = 2013
$ws.Range("F1").Value = 2014

# 3) The insert operation copied column A's style onto the new blank cells
#    in rows 2:9. The new data for those rows should be unstyled (matching
#    how the original B2:G9 cells looked), so clear that inherited format.
$ws.Range("B2:F9").ClearFormats()

# The old B column (now shifted to G) held an empty placeholder string for
# rows 3:9 (mean/std/min/25%/50%/75%/max have no value for "count"-less
# years). That placeholder is now stale (G has its own fresh "no data"
# entry like the rest of the new columns), so clear it out too.
$ws.Range("G3:G9").ClearContents()

# 4) New "count" row values (row 2) are all 0 for the newly added years.
$ws.Range("B2:G2").Value = 0

# 5) Rows 3-9 (mean/std/min/25%/50%/75%/max) have no data for the new
#    years - leave B3:G9 blank (same as the rest of that empty row).

# 6) A handful of the shifted statistics values differ slightly from the
#    old computation (recomputed stats), fix those explicit cells up.
$ws.Range("L2").Value = 43

$ws.Range("K3").Value = 5.82656053680151
$ws.Range("L3").Value = 5.565534543592364

$ws.Range("K4").Value = 0.9344216723656189
$ws.Range("L4").Value = 0.8958563041025963

$ws.Range("L6").Value = 4.774488785043467

$ws.Range("K7").Value = 5.816603052797736
$ws.Range("L7").Value = 5.468483350024568

$ws.Range("L8").Value = 6.297854220657186
